$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 74
$ws.Range("B7").Value = 74
$ws.Range("B17").Value = 74

$ws.Range("B17").Select()
